# Rename existing sheet "testing" -> "validScenario"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "validScenario"

# Duplicate sheet1 (carries over column widths, styles, shared header/data
# layout) right after itself, then rename it to the new scenario sheet.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "invalidScenario"

# Sheet2 (invalidScenario) keeps the same headers but an "invalid" data row;
# force the username as text so it matches the string-typed column.
$ws2.Range("A2").Value = "'1234"
[void]$ws2.Range("A2").Select()

# Sheet1 (validScenario) stays the active tab; only its selection moves
[void]$ws1.Activate()
[void]$ws1.Range("E10").Select()
